$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 2504
$ws.Range("C2").Value = 2471
$ws.Range("D2").Value = 2489.45
$ws.Range("E2").Value = 2496.5
$ws.Range("F2").Value = 74
$ws.Range("G2").Value = 2478.55

# Row 3
$ws.Range("B3").Value = 384.85
$ws.Range("C3").Value = 376
$ws.Range("D3").Value = 382
$ws.Range("E3").Value = 382.25
$ws.Range("F3").Value = 31
$ws.Range("G3").Value = 376.45

# Row 4
$ws.Range("B4").Value = 1522.45
$ws.Range("C4").Value = 1504.2
$ws.Range("D4").Value = 1519
$ws.Range("E4").Value = 1518.2
$ws.Range("F4").Value = 7
$ws.Range("G4").Value = 1510.45

# Row 5
$ws.Range("B5").Value = 7397
$ws.Range("C5").Value = 7296
$ws.Range("D5").Value = 7391
$ws.Range("E5").Value = 7387.8
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 7303.7

# Row 6
$ws.Range("B6").Value = 236.15
$ws.Range("C6").Value = 232.3
$ws.Range("D6").Value = 233.9
$ws.Range("E6").Value = 234
$ws.Range("F6").Value = 74
$ws.Range("G6").Value = 235.05

# Row 7
$ws.Range("B7").Value = 197.55
$ws.Range("C7").Value = 194
$ws.Range("D7").Value = 194.55
$ws.Range("E7").Value = 195
$ws.Range("F7").Value = 156
$ws.Range("G7").Value = 195.55

# Row 8
$ws.Range("B8").Value = 257.8
$ws.Range("C8").Value = 252.55
$ws.Range("D8").Value = 255.25
$ws.Range("E8").Value = 255.8
$ws.Range("F8").Value = 133
$ws.Range("G8").Value = 257.25

# Row 9
$ws.Range("B9").Value = 522.5
$ws.Range("C9").Value = 513.8
$ws.Range("D9").Value = 519
$ws.Range("E9").Value = 519
$ws.Range("F9").Value = 21
$ws.Range("G9").Value = 519.65

# Row 10
$ws.Range("B10").Value = 3417.1
$ws.Range("C10").Value = 3375
$ws.Range("D10").Value = 3390
$ws.Range("E10").Value = 3391.3
$ws.Range("F10").Value = 4
$ws.Range("G10").Value = 3408.4

# Row 11
$ws.Range("B11").Value = 145.5
$ws.Range("C11").Value = 142.85
$ws.Range("D11").Value = 143.9
$ws.Range("E11").Value = 143.8
$ws.Range("F11").Value = 173
$ws.Range("G11").Value = 144.7

# Row 12
$ws.Range("B12").Value = 1241.55
$ws.Range("C12").Value = 1222.7
$ws.Range("D12").Value = 1236
$ws.Range("E12").Value = 1235.85
$ws.Range("F12").Value = 17
$ws.Range("G12").Value = 1239.1

# Row 13
$ws.Range("B13").Value = 1600
$ws.Range("C13").Value = 1580.35
$ws.Range("D13").Value = 1599
$ws.Range("E13").Value = 1595.5
$ws.Range("F13").Value = 281
$ws.Range("G13").Value = 1582.25

# Row 14
$ws.Range("B14").Value = 480.55
$ws.Range("C14").Value = 474.25
$ws.Range("D14").Value = 478.1
$ws.Range("E14").Value = 478.15
$ws.Range("F14").Value = 51
$ws.Range("G14").Value = 479.55

# Row 15
$ws.Range("B15").Value = 966.35
$ws.Range("C15").Value = 950.8
$ws.Range("D15").Value = 960.05
$ws.Range("E15").Value = 957.5
$ws.Range("F15").Value = 308
$ws.Range("G15").Value = 964.45

# Row 16
$ws.Range("B16").Value = 1416.9
$ws.Range("C16").Value = 1396.45
$ws.Range("D16").Value = 1405.25
$ws.Range("E16").Value = 1408.65
$ws.Range("F16").Value = 40
$ws.Range("G16").Value = 1409.9

# Row 17
$ws.Range("B17").Value = 1480.25
$ws.Range("C17").Value = 1463.3
$ws.Range("D17").Value = 1475.8
$ws.Range("E17").Value = 1477.45
$ws.Range("F17").Value = 50
$ws.Range("G17").Value = 1475.05

# Row 18
$ws.Range("B18").Value = 704.45
$ws.Range("C18").Value = 697
$ws.Range("D18").Value = 701.25
$ws.Range("E18").Value = 702.45
$ws.Range("F18").Value = 16
$ws.Range("G18").Value = 701.55

# Row 19
$ws.Range("B19").Value = 454
$ws.Range("C19").Value = 445
$ws.Range("D19").Value = 449.6
$ws.Range("E19").Value = 449.25
$ws.Range("F19").Value = 28
$ws.Range("G19").Value = 452.8

# Row 20
$ws.Range("B20").Value = 1580.9
$ws.Range("C20").Value = 1558.7
$ws.Range("D20").Value = 1572.2
$ws.Range("E20").Value = 1577
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 1567.9

# Row 21
$ws.Range("B21").Value = 294.8
$ws.Range("C21").Value = 289.05
$ws.Range("D21").Value = 294
$ws.Range("E21").Value = 294.35
$ws.Range("F21").Value = 26
$ws.Range("G21").Value = 293.6

# Row 22
$ws.Range("B22").Value = 2436.15
$ws.Range("C22").Value = 2417.25
$ws.Range("D22").Value = 2426.35
$ws.Range("E22").Value = 2428.7
$ws.Range("F22").Value = 41
$ws.Range("G22").Value = 2428.45

# Row 23
$ws.Range("B23").Value = 573.85
$ws.Range("C23").Value = 568
$ws.Range("D23").Value = 570.95
$ws.Range("E23").Value = 570.5
$ws.Range("F23").Value = 206
$ws.Range("G23").Value = 572.3

# Row 24
$ws.Range("B24").Value = 619.55
$ws.Range("C24").Value = 607
$ws.Range("D24").Value = 614
$ws.Range("E24").Value = 615.35
$ws.Range("F24").Value = 9
$ws.Range("G24").Value = 614.15

# Row 25
$ws.Range("B25").Value = 1099.25
$ws.Range("C25").Value = 1078.2
$ws.Range("D25").Value = 1090
$ws.Range("E25").Value = 1090.1
$ws.Range("F25").Value = 6
$ws.Range("G25").Value = 1097.45

# Row 26
$ws.Range("B26").Value = 611.45
$ws.Range("C26").Value = 605.6
$ws.Range("D26").Value = 610.6
$ws.Range("E26").Value = 610.25
$ws.Range("F26").Value = 65
$ws.Range("G26").Value = 608.9

# Row 27
$ws.Range("B27").Value = 259.75
$ws.Range("C27").Value = 254.65
$ws.Range("D27").Value = 257.1
$ws.Range("E27").Value = 257.5
$ws.Range("F27").Value = 92
$ws.Range("G27").Value = 258.85

# Row 28
$ws.Range("B28").Value = 131.05
$ws.Range("C28").Value = 128.2
$ws.Range("D28").Value = 129.45
$ws.Range("E28").Value = 129.6
$ws.Range("F28").Value = 360
$ws.Range("G28").Value = 130.8

# Row 29
$ws.Range("B29").Value = 8556.799999999999
$ws.Range("C29").Value = 8444.1
$ws.Range("D29").Value = 8550
$ws.Range("E29").Value = 8537.4
$ws.Range("F29").Value = 3
$ws.Range("G29").Value = 8508.85
